$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 291, shifting existing rows 291-310 down to 292-311.
$ws.Rows.Item(291).Insert()

# Populate the newly inserted row 291 with the new record.
$ws.Cells.Item(291, 1).Value = 10
$ws.Cells.Item(291, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(291, 3).Value = "La Araucanía"
$ws.Cells.Item(291, 4).Value = 44826
$ws.Cells.Item(291, 5).Value = 9
$ws.Cells.Item(291, 6).Value = "Fruta"
$ws.Cells.Item(291, 7).Value = 100102
$ws.Cells.Item(291, 8).Value = "Cítricos"
$ws.Cells.Item(291, 9).Value = 100102006
$ws.Cells.Item(291, 10).Value = "Pomelo"
$ws.Cells.Item(291, 11).Value = "Start Ruby"
$ws.Cells.Item(291, 12).Value = "Primera"
$ws.Cells.Item(291, 13).Value = 95
$ws.Cells.Item(291, 14).Value = 12000
$ws.Cells.Item(291, 15).Value = 12000
$ws.Cells.Item(291, 16).Value = 12000
$ws.Cells.Item(291, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(291, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(291, 19).Value = 800
$ws.Cells.Item(291, 20).Value = 15
